$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @("LaMelo Ball", "PG,SG", "Charlotte Hornets")
    3  = @("Derrick White", "PG,SG", "Boston Celtics")
    6  = @("Brandon Miller", "SG,SF", "Charlotte Hornets")
    7  = @("Grant Williams", "PF,C", "Charlotte Hornets")
    8  = @("Anthony Davis", "PF,C", "Los Angeles Lakers")
    9  = @("Yves Missi", "C", "New Orleans Pelicans")
    10 = @("Moussa Diabate", "C", "Charlotte Hornets")
    11 = @("Julius Randle", "PF", "Minnesota Timberwolves")
    13 = @("Cameron Johnson", "SF,PF", "Brooklyn Nets")
    14 = @("Cade Cunningham", "PG,SG", "Detroit Pistons")
    15 = @("Cam Thomas", "SG,SF", "Brooklyn Nets")
    16 = @("Royce O'Neale", "SF,PF", "Phoenix Suns")
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
}
